$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Projects Final Score" in C1, bold font, no wrap/centering
$ws.Range("C1").Value = "Projects Final Score"
$ws.Range("C1").Font.Bold = $true

# Set column C width to match target (stored width ~18.6 chars wide,
# closest reachable value given Excel's column-width rounding granularity)
$ws.Columns.Item(3).ColumnWidth = 17.75

# Drop NaN rows -- renumber column A values (16,17,18,19,20) -> (11,12,13,14,15)
$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 13
$ws.Range("A5").Value = 14
$ws.Range("A6").Value = 15

# Move selection
$ws.Range("F10").Select()
